# Auto-generated: applies scheduled-runner value updates to Atomos_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 168663.17
$ws.Range("J45").Value = 168663.17
$ws.Range("L45").Value = 505989.51
$ws.Range("N45").Value = -506373.51

$ws.Range("H88").Value = 4953.52
$ws.Range("I88").Value = 4464.3076
$ws.Range("J88").Value = 5483.5
$ws.Range("K88").Value = 4464.3076
$ws.Range("L88").Value = 5483.5
$ws.Range("M88").Value = -4058.3076
$ws.Range("N88").Value = -6295.5

$ws.Range("H91").Value = 4953.52
$ws.Range("I91").Value = 4464.3076
$ws.Range("J91").Value = 5483.5
$ws.Range("K91").Value = 4464.3076
$ws.Range("L91").Value = 5483.5
$ws.Range("M91").Value = -3060.3076
$ws.Range("N91").Value = -8291.5

$ws.Range("H132").Value = 6455052.5
$ws.Range("I132").Value = 8003397
$ws.Range("J132").Value = 3615.8333
$ws.Range("K132").Value = 24010191
$ws.Range("L132").Value = 10847.4999
$ws.Range("M132").Value = -24007661
$ws.Range("N132").Value = -15907.4999

$ws.Range("H138").Value = 5039.5874
$ws.Range("I138").Value = 2334.1072
$ws.Range("J138").Value = 6496.385
$ws.Range("K138").Value = 7002.321599999999
$ws.Range("L138").Value = 19489.155
$ws.Range("M138").Value = -1862.321599999999
$ws.Range("N138").Value = -29769.155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 150005000
$ws.Range("I2").Value = 250002340
$ws.Range("K2").Value = 250002340
$ws.Range("M2").Value = -250002227

$ws.Range("H23").Value = 62879
$ws.Range("J23").Value = 52602.8
$ws.Range("L23").Value = 52602.8
$ws.Range("N23").Value = -53120.8

$ws.Range("H44").Value = 18839.6
$ws.Range("J44").Value = 21924.5
$ws.Range("L44").Value = 21924.5
$ws.Range("N44").Value = -22900.5

$ws.Range("H55").Value = 15050.5
$ws.Range("I55").Value = 2048
$ws.Range("K55").Value = 2048
$ws.Range("M55").Value = -1733

$ws.Range("H74").Value = 2611.5
$ws.Range("I74").Value = 1827.5
$ws.Range("J74").Value = 5747.5
$ws.Range("K74").Value = 1827.5
$ws.Range("L74").Value = 5747.5
$ws.Range("M74").Value = -953.5
$ws.Range("N74").Value = -7495.5

$ws.Range("H77").Value = 2611.5
$ws.Range("I77").Value = 1827.5
$ws.Range("J77").Value = 5747.5
$ws.Range("K77").Value = 9137.5
$ws.Range("L77").Value = 28737.5
$ws.Range("M77").Value = -4769.5
$ws.Range("N77").Value = -37473.5

$ws.Range("H116").Value = 150005000
$ws.Range("I116").Value = 250002340
$ws.Range("K116").Value = 250002340
$ws.Range("M116").Value = -250000046

$ws.Range("H122").Value = 3323.5
$ws.Range("I122").Value = 2294.077
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 6882.231000000001
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -4432.231000000001
$ws.Range("N122").Value = -22900

$ws.Range("H132").Value = 1891.75
$ws.Range("I132").Value = 1458.84
$ws.Range("J132").Value = 5499.3335
$ws.Range("K132").Value = 4376.52
$ws.Range("L132").Value = 16498.0005
$ws.Range("M132").Value = -1846.52
$ws.Range("N132").Value = -21558.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 150005000
$ws.Range("I3").Value = 250002340
$ws.Range("K3").Value = 250002340
$ws.Range("M3").Value = -250002226

$ws.Range("H94").Value = 1095.4348
$ws.Range("I94").Value = 637.375
$ws.Range("J94").Value = 2142.4285
$ws.Range("K94").Value = 637.375
$ws.Range("L94").Value = 2142.4285
$ws.Range("M94").Value = -186.375
$ws.Range("N94").Value = -3044.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 45021
$ws.Range("J29").Value = 45021
$ws.Range("L29").Value = 45021
$ws.Range("N29").Value = -45607

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 2633
$ws.Range("J93").Value = 3000
$ws.Range("L93").Value = 9000
$ws.Range("N93").Value = -12744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 37842.855
$ws.Range("J141").Value = 37842.855
$ws.Range("L141").Value = 37842.855
$ws.Range("N141").Value = -48202.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 10000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -10224

$ws.Range("H22").Value = 1354.7273
$ws.Range("I22").Value = 425
$ws.Range("J22").Value = 1886
$ws.Range("K22").Value = 425
$ws.Range("L22").Value = 1886
$ws.Range("M22").Value = -130
$ws.Range("N22").Value = -2476

$ws.Range("H27").Value = 1354.7273
$ws.Range("I27").Value = 425
$ws.Range("J27").Value = 1886
$ws.Range("K27").Value = 425
$ws.Range("L27").Value = 1886
$ws.Range("M27").Value = -318
$ws.Range("N27").Value = -2100

$ws.Range("H40").Value = 8349.083000000001
$ws.Range("I40").Value = 8518.9
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 8518.9
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -8382.9
$ws.Range("N40").Value = -7772

$ws.Range("H46").Value = 2837.5
$ws.Range("J46").Value = 4300
$ws.Range("L46").Value = 4300
$ws.Range("N46").Value = -4676

$ws.Range("H122").Value = 3639.261
$ws.Range("I122").Value = 2780.8
$ws.Range("J122").Value = 4299.615
$ws.Range("K122").Value = 8342.400000000001
$ws.Range("L122").Value = 12898.845
$ws.Range("M122").Value = -5892.400000000001
$ws.Range("N122").Value = -17798.845

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 30000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -34940

$ws.Range("H127").Value = 29663.572
$ws.Range("J127").Value = 29663.572
$ws.Range("L127").Value = 29663.572
$ws.Range("N127").Value = -39583.572

$ws.Range("H136").Value = 3227.1135
$ws.Range("I136").Value = 2231.8708
$ws.Range("J136").Value = 5600.385
$ws.Range("K136").Value = 6695.6124
$ws.Range("L136").Value = 16801.155
$ws.Range("M136").Value = -4145.6124
$ws.Range("N136").Value = -21901.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 969.8333
$ws.Range("I113").Value = 129.66667
$ws.Range("J113").Value = 1810
$ws.Range("K113").Value = 389.00001
$ws.Range("L113").Value = 5430
$ws.Range("M113").Value = 1780.99999
$ws.Range("N113").Value = -9770

$ws.Range("H126").Value = 47115.305
$ws.Range("I126").Value = 74333.28999999999
$ws.Range("J126").Value = 4776.222
$ws.Range("K126").Value = 222999.87
$ws.Range("L126").Value = 14328.666
$ws.Range("M126").Value = -220529.87
$ws.Range("N126").Value = -19268.666

$ws.Range("H141").Value = 27892.857
$ws.Range("J141").Value = 27892.857
$ws.Range("L141").Value = 27892.857
$ws.Range("N141").Value = -38252.857
